$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A17").Value = "DB auto backups"
$ws.Range("B17").Value = "Ignas"
$ws.Range("C17").Value = "vidutinis"

$ws.Range("A18").Value = "Taisyklės"
$ws.Range("B18").Value = "Ignas"
$ws.Range("C18").Value = "vidutinis"

$ws.Range("A23").Select()
